$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 - Bolivia
$ws.Range("B53").Value = 16165
$ws.Range("C53").Value = 884
$ws.Range("D53").Value = 2372
$ws.Range("E53").Value = 13260
$ws.Range("G53").Value = 21
$ws.Range("H53").Value = 533

# Row 70 - Honduras
$ws.Range("B70").Value = 7669
$ws.Range("C70").Value = 309
$ws.Range("D70").Value = 837
$ws.Range("E70").Value = 6538
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 294

# Row 71 - Australia
$ws.Range("B71").Value = 7288
$ws.Range("C71").Value = 3
$ws.Range("D71").Value = 6777
$ws.Range("E71").Value = 409

# Row 87 - El Salvador
$ws.Range("D87").Value = 1504
$ws.Range("E87").Value = 1803
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 66
